$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Orders")

# Row 2: Chai order - note text renamed to new order number
$ws.Range("D2").Value = "Order #512711"

# Row 3: Ipoh Coffee - was "Unsuccessful / out of stock", now Successful with an order
# number and quantities filled in
$ws.Range("C3").Value = "Successful"
$ws.Range("D3").Value = "Order #512711"
$ws.Range("E3").Value = 23
$ws.Range("F3").Value = 69

# Rows 4-8: same order-number rename
$ws.Range("D4").Value = "Order #512711"
$ws.Range("D5").Value = "Order #512711"
$ws.Range("D6").Value = "Order #512711"
$ws.Range("D7").Value = "Order #512711"
$ws.Range("D8").Value = "Order #512711"

# New cell H2 holding a total, formatted like the existing currency cells
$ws.Range("H2").Value = 233
$ws.Range("H2").NumberFormat = $ws.Range("F2").NumberFormat

# Rows 10 and 11 lose their Status/Notes/Unit Price/Total Price data, becoming
# blank placeholder rows like rows 12-13
$ws.Range("C10:D11").ClearContents()
$ws.Range("E10:F11").ClearContents()
